$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "results-1606079248"

# Add new header columns
$ws.Range("F1").Value = "cplex time (sec)"
$ws.Range("G1").Value = "cplex time (ms)"
$ws.Range("H1").Value = "result"

# Update the per-graph data rows: D (heuristic_result dev?), E (heuristic_time ms),
# F (cplex time sec, new), G (cplex time ms, new), H (result, new)
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 29
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 14
$ws.Range("H2").Value = 32
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 12
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 25
$ws.Range("H3").Value = 16
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 73
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 2222
$ws.Range("H4").Value = 12
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 82
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1125
$ws.Range("H5").Value = 24
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 361
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 185
$ws.Range("H6").Value = 55
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 90
$ws.Range("H7").Value = 4
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 128
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 116
$ws.Range("H8").Value = 30
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 506
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 43
$ws.Range("H9").Value = 70
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 366
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 54
$ws.Range("H10").Value = 60
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1915
$ws.Range("F11").Value = 4
$ws.Range("G11").Value = 4429
$ws.Range("H11").Value = 126
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 744
$ws.Range("F12").Value = 7
$ws.Range("G12").Value = 7155
$ws.Range("H12").Value = 64
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 211
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1926
$ws.Range("H13").Value = 58
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 407
$ws.Range("F14").Value = 12
$ws.Range("G14").Value = 12560
$ws.Range("H14").Value = 14
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 466
$ws.Range("F15").Value = 14
$ws.Range("G15").Value = 14702
$ws.Range("H15").Value = 26
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 122
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 317
$ws.Range("H16").Value = 18
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 118
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 470
$ws.Range("H17").Value = 34
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 61
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 1951
$ws.Range("H18").Value = 11
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 76
$ws.Range("F19").Value = 15
$ws.Range("G19").Value = 15854
$ws.Range("H19").Value = 12
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 161
$ws.Range("F20").Value = 25
$ws.Range("G20").Value = 25443
$ws.Range("H20").Value = 8
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 92
$ws.Range("F21").Value = 28
$ws.Range("G21").Value = 28362
$ws.Range("H21").Value = 15
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 374
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 771
$ws.Range("H22").Value = 44
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 111
$ws.Range("F23").Value = 35
$ws.Range("G23").Value = 35972
$ws.Range("H23").Value = 17
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 327
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 136
$ws.Range("H24").Value = 44
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 145
$ws.Range("F25").Value = 164
$ws.Range("G25").Value = 164488
$ws.Range("H25").Value = 21
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 316
$ws.Range("F26").Value = 48
$ws.Range("G26").Value = 48240
$ws.Range("H26").Value = 25
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 138
$ws.Range("F27").Value = 64
$ws.Range("G27").Value = 64770
$ws.Range("H27").Value = 18

# Update the active selection as recorded in the saved view
$ws.Range("L8").Select()
